$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.053.64"
$ws.Range("E2").Value = "  +2.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.787.60"
$ws.Range("E3").Value = "  +2.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +1.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.55"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.561"
$ws.Range("E6").Value = "  +3.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +1.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.20"
$ws.Range("E8").Value = "  -5.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.50"
$ws.Range("E9").Value = "  +2.43%  "

$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.043.87"
$ws.Range("E13").Value = "  +3.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.782.76"
$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.628"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.070.87"
$ws.Range("E16").Value = "  +3.37%  "

$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.19"
$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.99"
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "252.55"
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0742"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("E22").Value = "  +2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.37"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("E24").Value = "  -2.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("E25").Value = "  -1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.48"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.56"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0515"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("E35").Value = "  +4.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.506.95"
$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("E37").Value = "  +2.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.637"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.49"
$ws.Range("E40").Value = "  -1.77%  "

$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.902"
$ws.Range("E43").Value = "  +4.11%  "

$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.07"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.941.46"
$ws.Range("E47").Value = "  +3.42%  "

$ws.Range("E48").Value = "  +2.40%  "

$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.86"
$ws.Range("E50").Value = "  +13.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.50"
$ws.Range("E51").Value = "  -6.88%  "
